# Outline.docx edit script
# Summary of content changes (per the commit's xml diff):
#  1. Paragraph 2: "fine grained" -> "high quality"
#     (the Word "_GoBack" bookmark also relocates to sit mid-word in
#     "complex", reflecting where the author's cursor ended up editing)
#  2. Paragraph 4: trailing two spaces after "...annotators, and " get
#     replaced with a new sentence about engaging less experienced
#     observers.
#  3. A brand new paragraph is inserted right after that one: "While
#     computer vision has made incredible strides ... Computer vision
#     articles  " (trailing run is a single extra space).
# (The remaining hunks in the diff only wrap existing words such as
#  "TensorFlow", "tensorboard", "docker", "cloudml" in spell/grammar-
#  check <w:proofErr/> markers - no visible text changes there, so
#  there is nothing further to type.)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "fine grained" -> "high quality" in paragraph 2, and move the
#    _GoBack bookmark to land between "co" and "mplex" (matching the
#    diff's run split), same as it would after a live edit session.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("fine grained", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
if (-not $rng.Find.Found) {
    throw "Could not find 'fine grained' to replace"
}
$start = $rng.Start
$rng.Text = ""
$insPt = $d.Range($start, $start)
$insPt.InsertAfter("high quality")

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$complexRng = $d.Content
$complexRng.Find.Execute("complex logistics", $true, $false, $false, $false, `
                          $false, $true, 1, $false, "", 0)
if (-not $complexRng.Find.Found) {
    throw "Could not find 'complex logistics' to re-anchor bookmark"
}
$bmPoint = $d.Range($complexRng.Start + 2, $complexRng.Start + 2)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# ---------------------------------------------------------------------
# 2. Extend paragraph 4: replace the trailing "  " (two spaces) after
#    "...annotators, and " with the new sentence.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4Rng = $p4.Range
$p4Rng.Find.Execute("and   ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, `
                     "and engage less experienced observers in biodiversity monitoring. Computer vision is a field of image-based computer science that uses image pixels to mimic human perception based on image characteristics, shape and sequence. ", `
                     2)
if (-not $p4Rng.Find.Found) {
    throw "Could not find paragraph 4 trailing spaces to replace"
}

# ---------------------------------------------------------------------
# 3. Insert a brand new paragraph right after (now) paragraph 4 with
#    the "While computer vision..." text.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$endOfP4 = $p4.Range
$endOfP4.Collapse(0)
$endOfP4.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.InsertAfter("While computer vision has made incredible strides in a variety of fields, its growth in ecology has been slowed by a lack of access to high level algorithms. Computer vision articles  ")

Write-Host "Paragraph 2: [$($d.Paragraphs(2).Range.Text)]"
Write-Host "Paragraph 4: [$($d.Paragraphs(4).Range.Text)]"
Write-Host "Paragraph 5: [$($d.Paragraphs(5).Range.Text)]"
Write-Host "Total paragraphs: $($d.Paragraphs.Count)"
